$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J1").Value = "Git Commit"
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 10).Value = "2b20aa4"
}

$ws.Range("J17").Select()
